$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.575.16"
$ws.Range("D3").Value = "2.288.64"
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "95.38"
$ws.Range("E5").Value = "  -3.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "267.87"
$ws.Range("E6").Value = "  -2.54%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.623"
$ws.Range("E7").Value = "  -0.81%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("E8").Value = "  -0.10%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.609"
$ws.Range("E9").Value = "  -4.05%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "44.62"
$ws.Range("E10").Value = "  -7.62%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0938"
$ws.Range("E11").Value = "  -1.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.80"
$ws.Range("E12").Value = "  -5.22%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.104"
$ws.Range("E13").Value = "  -0.43%  "
$ws.Range("D14").Value = "2.632.12"
$ws.Range("E14").Value = "  +1.35%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.18"
$ws.Range("E15").Value = "  -2.57%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.852"
$ws.Range("E16").Value = "  +1.33%  "
$ws.Range("D17").Value = "2.290.87"
$ws.Range("E17").Value = "  +0.68%  "
$ws.Range("D18").Value = "43.548.47"
$ws.Range("E18").Value = "  -1.47%  "
$ws.Range("E19").Value = "  +0.38%  "
$ws.Range("E20").Value = "  -0.48%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "72.46"
$ws.Range("E21").Value = "  +2.07%  "
$ws.Range("E22").Value = "  +3.72%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "235.15"
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.04"
$ws.Range("E24").Value = "  -17.13%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("E26").Value = "  -1.22%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.22"
$ws.Range("E27").Value = "  -2.83%  "
$ws.Range("E28").Value = "  +1.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.03"
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("E30").Value = "  -0.04%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "175.18"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.99"
$ws.Range("E32").Value = "  +3.41%  "
$ws.Range("E33").Value = "  -4.22%  "
$ws.Range("E34").Value = "  -7.17%  "
$ws.Range("E35").Value = "  -0.20%  "
$ws.Range("E36").Value = "  +0.55%  "
$ws.Range("E37").Value = "  -5.77%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.39"
$ws.Range("E38").Value = "  -1.15%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.31"
$ws.Range("E39").Value = "  -8.29%  "
$ws.Range("B40").Value = "Algorand"
$ws.Range("C40").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.236"
$ws.Range("E40").Value = "  -7.41%  "
$ws.Range("B41").Value = "LidoDAOToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.34"
$ws.Range("E41").Value = "  +6.15%  "
$ws.Range("E42").Value = "  +15.08%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "64.57"
$ws.Range("E43").Value = "  +3.02%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "12.03"
$ws.Range("E44").Value = "  -5.22%  "
$ws.Range("E45").Value = "  +2.71%  "
$ws.Range("E46").Value = "  -4.81%  "
$ws.Range("E47").Value = "  -1.40%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "98.18"
$ws.Range("E48").Value = "  -2.31%  "
$ws.Range("E49").Value = "  -0.75%  "
$ws.Range("B50").Value = "RocketPoolETH"
$ws.Range("C50").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D50").Value = "2.511.73"
$ws.Range("E50").Value = "  +1.34%  "
$ws.Range("B51").Value = "WOONetwork"
$ws.Range("C51").Value = "https://coinranking.com/coin/k-J3YwacF+woonetwork-woo"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.429"
$ws.Range("E51").Value = "  -2.14%  "
